$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the Stewart et al. (2021) albula/lavaretus rows (5 & 6) ---
$ws.Range("C5").Value = "albula"
$ws.Range("F5").Value = -2.366372
$ws.Range("G5").Value = 0.008832
$ws.Range("H5").Value = 0.005009

$ws.Range("C6").Value = "lavaretus"
$ws.Range("F6").Value = -2.418263
$ws.Range("G6").Value = 0.045934
$ws.Range("H6").Value = 0.003216

# --- Swap the Eckmann (1987) Lake Geneva/Lake Bourget rows (12 & 13) ---
$ws.Range("E12").Value = "Lake Geneva"
$ws.Range("F12").Value = -2.115944
$ws.Range("G12").Value = 0.052842

$ws.Range("E13").Value = "Lake Bourget"
$ws.Range("F13").Value = -2.06886
$ws.Range("G13").Value = 0.047427

# --- Update the current selection ---
$ws.Range("G17").Select()

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1
